$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

function Get-ParagraphRangeByText($needle) {
    $paras = $d.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        $rr = $p.Range
        if ($rr.Text -like "*$needle*") {
            return $rr
        }
    }
    return $null
}

# Apply a split: locate the paragraph containing $origRunText (that text must be
# the literal, currently-existing substring produced by a single run inside that
# paragraph), replace from the start of that text through the end of the
# paragraph with $newRunsXml, followed by whatever XML already followed that
# text in the paragraph (so trailing runs/formatting are preserved untouched
# and are not reordered).
function Apply-RunSplit($origRunText, $newRunsXml) {
    $p = Get-ParagraphRangeByText($origRunText)
    $fullText = $p.Text
    $off = $fullText.IndexOf($origRunText)
    $runStart = $p.Start + $off
    $runEnd = $p.End
    $target = $d.Range($runStart, $runEnd)

    $escaped = [regex]::Escape($origRunText)
    $docXml = $d.Content.WordOpenXML
    $pattern = [regex]::Escape("<w:t>") + $escaped + [regex]::Escape("</w:t></w:r>") + "(.*?)</w:p>"
    $trailing = ""
    if ($docXml -match $pattern) {
        $trailing = $matches[1]
    } else {
        $pattern2 = [regex]::Escape("<w:t xml:space=`"preserve`">") + $escaped + [regex]::Escape("</w:t></w:r>") + "(.*?)</w:p>"
        if ($docXml -match $pattern2) {
            $trailing = $matches[1]
        }
    }

    $xml = '<w:p xmlns:w="' + $wNs + '">' + $newRunsXml + $trailing + '</w:p>'
    $target.InsertXML($xml)
}

# --- Edit 1: "posto )" paragraph - split off "posto )" with gramStart/gramEnd proofErr ---
Apply-RunSplit " O Custos comtemplará tabelas congeladas/fixas, devido o foco da aplicação ser o controle de veículo, não será validado custos vinculados à externos (exemplo: posto )" (
    '<w:r w:rsidR="00E4386E"><w:t xml:space="preserve"> O Custos comtemplará tabelas congeladas/fixas, devido o foco da aplicação ser o controle de veículo, não será validado custos vinculados à externos (exemplo: </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r w:rsidR="00E4386E"><w:t>posto )</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>'
)

# --- Edit 2: "(Apresentar Romanhuk)" paragraph - split off "Romanhuk" with spellStart/spellEnd proofErr ---
Apply-RunSplit "(Apresentar Romanhuk)" (
    '<w:r><w:t xml:space="preserve">(Apresentar </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Romanhuk</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>)</w:t></w:r>'
)

# --- Edit 3: "A partir de: 26/02 codar." paragraph - split off "codar" with spellStart/spellEnd proofErr ---
Apply-RunSplit "A partir de: 26/02 codar." (
    '<w:r><w:t xml:space="preserve">A partir de: 26/02 </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>codar</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>.</w:t></w:r>'
)

Write-Output "All edits applied"
